# 55927 - "handle date formulas too"
#
# The original sheet has a single XML-mapped cell, A1, holding a date
# value. This adds a second row (A2) containing a DATE() formula whose
# result also needs to flow through the same "date" formatting as A1,
# mirroring how Marc's updated patch taught the XML export path to cope
# with formula-driven dates as well as literal ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone A1's cell format (built-in short-date number format) onto A2 via
# copy/paste-special so we reuse the existing style record instead of
# minting a new one, then drop in the date formula itself.
$ws.Range("A1").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A2").Formula = "=DATE(2012,2,16)"

# Match the post-entry selection Excel leaves behind: the cursor advances
# to the next empty row after typing the formula into A2.
$ws.Range("A3").Select()
